$wb = $excel.ActiveWorkbook

# --- Update the conversion note text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.73 = 51438.88 pesos`n✅ 51438.88 pesos = 12.72 = 982.77 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate cells on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 78.54000000000001
$wsTasas.Range("O10").Value = 4040.01
$wsTasas.Range("N12").Value = 4045
$wsTasas.Range("O12").Value = 77.282
